# Issue on reindexing when generating cvxpy variables solved (see issue #23 on github)
#
# This script reproduces, via the Excel COM object model, the edits made to
# default/2_sut_multi_year_rcot/concept.xlsx, sheet "multi-year":
#   - the "I_ff" label and its 1,1 matrix (G7 / G8:H8) are relocated to new
#     rows 34/35 (G34 / G35:H35), fixing a reindexing bug
#   - the S8 array formula is repointed to the relocated matrix
#   - the N7 growth-factor cell and the N8:Q8 shared formula that depended on
#     it are removed, the resulting values are hard-coded instead
#   - a couple of hand-entered numbers are corrected (rounded) which cascades
#     through the dependent array formulas
#   - a blank styled spacer cell is added at L1
#   - the view is changed: gridlines hidden, zoomed to 130%, new selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("multi-year")

# ---------------------------------------------------------------------
# 1. Relocate the "I_ff" label (G7) and its 1/1 matrix (G8:H8) down to the
#    new rows 34/35, preserving their formatting.
# ---------------------------------------------------------------------
$labelText = $ws.Range("G7").Value2

$ws.Range("G8:H8").Copy()
$ws.Range("G35:H35").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("G34").Value2 = $labelText
$ws.Range("G35").Value2 = 1
$ws.Range("H35").Value2 = 1

# Point the S8:W8 array formula at the relocated matrix instead of G8:H8.
$ws.Range("S8:W8").FormulaArray = "=MMULT(G35:H35,TRANSPOSE(G16:H20))"

# Remove the now-vacated source cells.
$ws.Range("G7").Clear()
$ws.Range("G8:H8").Clear()

# ---------------------------------------------------------------------
# 2. Drop the growth-factor cell N7 and convert the N8:Q8 shared formula
#    (which depended on it) into plain cached numbers. Before clearing
#    N7, reuse its formatting for the new blank spacer cell at L1.
# ---------------------------------------------------------------------
$ws.Range("N7").Copy()
$ws.Range("L1").PasteSpecial(-4122)        # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("L1").ClearContents()

$ws.Range("N8").Value2 = 120
$ws.Range("O8").Value2 = 144
$ws.Range("P8").Value2 = 173
$ws.Range("Q8").Value2 = 207

$ws.Range("N7").Clear()

# ---------------------------------------------------------------------
# 3. Correct a couple of hand-entered figures; this cascades through the
#    dependent array formulas (S8:W8, J23:K27, S17) automatically.
# ---------------------------------------------------------------------
$ws.Range("V11").Value2 = 128
$ws.Range("W11").Value2 = 207

$ws.Range("H19").Value2 = 128
$ws.Range("H20").Value2 = 207

$wb.Application.CalculateFull()

# ---------------------------------------------------------------------
# 4. View changes: hide gridlines, zoom to 130%, move the selection.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $false
$win.Zoom = 130
$ws.Range("AE14").Select()

$wb.Application.CalculateFull()
